$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create new row 8 by duplicating row 7's formatting/layout (keeps B8/C8 blank
# cells present and copies A8's style from A7), then overwrite its contents.
$ws.Range("A7:D7").Copy($ws.Range("A8:D8"))

# Row 8 becomes "Other" with the new D8 value
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 344.8442797894652

# Row 7 becomes "Biogas" with its corrected D7 value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 357.6313357191123

# D3 is cleared (was a number, now blank)
$ws.Range("D3").Value = $null

# C4 and C5 corrected to 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
